$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 659.6
$ws.Range("J58").Value = 748
$ws.Range("L58").Value = 2244
$ws.Range("N58").Value = -2544
$ws.Range("H98").Value = 363.18182
$ws.Range("I98").Value = 366.875
$ws.Range("J98").Value = 353.33334
$ws.Range("K98").Value = 366.875
$ws.Range("L98").Value = 353.33334
$ws.Range("M98").Value = 1131.125
$ws.Range("N98").Value = -3349.33334
$ws.Range("H113").Value = 4639.9287
$ws.Range("I113").Value = 3476.3
$ws.Range("J113").Value = 7549
$ws.Range("K113").Value = 3476.3
$ws.Range("L113").Value = 7549
$ws.Range("M113").Value = -222.3000000000002
$ws.Range("N113").Value = -14057
$ws.Range("H122").Value = 363.18182
$ws.Range("I122").Value = 366.875
$ws.Range("J122").Value = 353.33334
$ws.Range("K122").Value = 1100.625
$ws.Range("L122").Value = 1060.00002
$ws.Range("M122").Value = 1349.375
$ws.Range("N122").Value = -5960.000019999999
$ws.Range("H132").Value = 5207.95
$ws.Range("I132").Value = 5231.278
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 15693.834
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -13163.834
$ws.Range("N132").Value = -20054
$ws.Range("H137").Value = 1706.9375
$ws.Range("I137").Value = 1630.5
$ws.Range("K137").Value = 4891.5
$ws.Range("M137").Value = -2341.5
$ws.Range("H138").Value = 1901.3673
$ws.Range("I138").Value = 2107.875
$ws.Range("J138").Value = 1861.0731
$ws.Range("K138").Value = 6323.625
$ws.Range("L138").Value = 5583.219300000001
$ws.Range("M138").Value = -1183.625
$ws.Range("N138").Value = -15863.2193

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7333
$ws.Range("I61").Value = 7333
$ws.Range("K61").Value = 7333
$ws.Range("M61").Value = -7121
$ws.Range("H63").Value = 3663.027
$ws.Range("I63").Value = 2323.7097
$ws.Range("J63").Value = 10582.833
$ws.Range("K63").Value = 2323.7097
$ws.Range("L63").Value = 10582.833
$ws.Range("M63").Value = -1637.7097
$ws.Range("N63").Value = -11954.833
$ws.Range("H66").Value = 3663.027
$ws.Range("I66").Value = 2323.7097
$ws.Range("J66").Value = 10582.833
$ws.Range("K66").Value = 11618.5485
$ws.Range("L66").Value = 52914.165
$ws.Range("M66").Value = -8186.548499999999
$ws.Range("N66").Value = -59778.165
$ws.Range("H95").Value = 54617
$ws.Range("J95").Value = 54617
$ws.Range("L95").Value = 54617
$ws.Range("N95").Value = -60109
$ws.Range("H110").Value = 4409.5386
$ws.Range("I110").Value = 4910.1816
$ws.Range("K110").Value = 4910.1816
$ws.Range("M110").Value = -2865.1816
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H136").Value = 7333
$ws.Range("I136").Value = 7333
$ws.Range("K136").Value = 21999
$ws.Range("M136").Value = -19449

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 69990
$ws.Range("J20").Value = 69990
$ws.Range("L20").Value = 69990
$ws.Range("N20").Value = -70462
$ws.Range("H30").Value = 69990
$ws.Range("J30").Value = 69990
$ws.Range("L30").Value = 69990
$ws.Range("N30").Value = -70172
$ws.Range("H58").Value = 3202.375
$ws.Range("I58").Value = 2677.25
$ws.Range("K58").Value = 2677.25
$ws.Range("M58").Value = -2474.25
$ws.Range("H99").Value = 5624.75
$ws.Range("I99").Value = 5624.75
$ws.Range("K99").Value = 5624.75
$ws.Range("M99").Value = -4126.75
$ws.Range("H126").Value = 5624.75
$ws.Range("I126").Value = 5624.75
$ws.Range("K126").Value = 16874.25
$ws.Range("M126").Value = -14404.25
$ws.Range("H128").Value = 69990
$ws.Range("J128").Value = 69990
$ws.Range("L128").Value = 69990
$ws.Range("N128").Value = -79950
$ws.Range("H129").Value = 94949
$ws.Range("J129").Value = 94949
$ws.Range("L129").Value = 94949
$ws.Range("N129").Value = -104949
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970
$ws.Range("H134").Value = 3412.3333
$ws.Range("I134").Value = 3457.625
$ws.Range("K134").Value = 10372.875
$ws.Range("M134").Value = -7837.875
$ws.Range("H136").Value = 3202.375
$ws.Range("I136").Value = 2677.25
$ws.Range("K136").Value = 8031.75
$ws.Range("M136").Value = -5481.75

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 575.86664
$ws.Range("I121").Value = 648.3333
$ws.Range("J121").Value = 467.16666
$ws.Range("K121").Value = 1944.9999
$ws.Range("L121").Value = 1401.49998
$ws.Range("M121").Value = -634.9999
$ws.Range("N121").Value = -4021.49998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4179
$ws.Range("I102").Value = 4490.2856
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 4490.2856
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -2868.2856
$ws.Range("N102").Value = -5244
$ws.Range("H126").Value = 2245.9
$ws.Range("I126").Value = 2245.9
$ws.Range("K126").Value = 6737.700000000001
$ws.Range("M126").Value = -4267.700000000001
$ws.Range("H132").Value = 3109.3333
$ws.Range("I132").Value = 3037.5652
$ws.Range("K132").Value = 9112.695599999999
$ws.Range("M132").Value = -6582.695599999999

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1871.7693
$ws.Range("I7").Value = 1613.3
$ws.Range("K7").Value = 1613.3
$ws.Range("M7").Value = -1501.3
$ws.Range("H40").Value = 6979.2666
$ws.Range("I40").Value = 4014.5386
$ws.Range("J40").Value = 26250
$ws.Range("K40").Value = 4014.5386
$ws.Range("L40").Value = 26250
$ws.Range("M40").Value = -3878.5386
$ws.Range("N40").Value = -26522
$ws.Range("H61").Value = 2923
$ws.Range("I61").Value = 2413.5
$ws.Range("K61").Value = 2413.5
$ws.Range("M61").Value = -2211.5
$ws.Range("H103").Value = 36601
$ws.Range("J103").Value = 36601
$ws.Range("L103").Value = 36601
$ws.Range("N103").Value = -38945
$ws.Range("H104").Value = 20542.5
$ws.Range("J104").Value = 20542.5
$ws.Range("L104").Value = 20542.5
$ws.Range("N104").Value = -27530.5
$ws.Range("H113").Value = 2923
$ws.Range("I113").Value = 2413.5
$ws.Range("K113").Value = 2413.5
$ws.Range("M113").Value = -243.5
$ws.Range("H126").Value = 1871.7693
$ws.Range("I126").Value = 1613.3
$ws.Range("K126").Value = 4839.9
$ws.Range("M126").Value = -2369.9
$ws.Range("H130").Value = 19109
$ws.Range("J130").Value = 19109
$ws.Range("L130").Value = 19109
$ws.Range("N130").Value = -29149
$ws.Range("H132").Value = 2221.3333
$ws.Range("I132").Value = 2237.5386
$ws.Range("J132").Value = 2179.2
$ws.Range("K132").Value = 6712.6158
$ws.Range("L132").Value = 6537.599999999999
$ws.Range("M132").Value = -4182.6158
$ws.Range("N132").Value = -11597.6
$ws.Range("H136").Value = 2624.8333
$ws.Range("I136").Value = 2749.9
$ws.Range("K136").Value = 8249.700000000001
$ws.Range("M136").Value = -5699.700000000001

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3825.3
$ws.Range("I62").Value = 3790.4
$ws.Range("J62").Value = 3860.2
$ws.Range("K62").Value = 3790.4
$ws.Range("L62").Value = 3860.2
$ws.Range("M62").Value = -3166.4
$ws.Range("N62").Value = -5108.2
$ws.Range("H65").Value = 3825.3
$ws.Range("I65").Value = 3790.4
$ws.Range("J65").Value = 3860.2
$ws.Range("K65").Value = 18952
$ws.Range("L65").Value = 19301
$ws.Range("M65").Value = -15832
$ws.Range("N65").Value = -25541
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H125").Value = 44375
$ws.Range("J125").Value = 44375
$ws.Range("L125").Value = 44375
$ws.Range("N125").Value = -54215
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
